$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) from column R into the new column S
# for every row that currently has a style applied there.
$ws.Range("R2:R6").Copy()
$ws.Range("S2:S6").PasteSpecial(-4122)

# Populate the new column's values (2022 figures)
$ws.Range("S3").Value = 2022
$ws.Range("S4").Value = 265803
$ws.Range("S5").Value = 3.8
$ws.Range("S6").Value = 33.6

# Update the selection to match the saved workbook state
$ws.Range("C19").Select()
